# Auto-generated PowerShell Excel COM-interop script
# Updates column F ("想去人数") values across sheets 展览, 演出, 全部类型

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 408   # was 407
$ws.Range("F5").Value = 5344   # was 5331
$ws.Range("F6").Value = 5344   # was 5331
$ws.Range("F7").Value = 176   # was 174
$ws.Range("F9").Value = 536   # was 535
$ws.Range("F11").Value = 1194   # was 1192
$ws.Range("F12").Value = 5291   # was 5286
$ws.Range("F15").Value = 102   # was 101
$ws.Range("F16").Value = 2950   # was 2883
$ws.Range("F20").Value = 3999   # was 3994
$ws.Range("F24").Value = 3916   # was 3912
$ws.Range("F29").Value = 254   # was 253
$ws.Range("F31").Value = 114   # was 113
$ws.Range("F32").Value = 114   # was 113
$ws.Range("F36").Value = 27   # was 26
$ws.Range("F37").Value = 6963   # was 6950
$ws.Range("F38").Value = 28   # was 26
$ws.Range("F39").Value = 1141   # was 1138
$ws.Range("F40").Value = 547   # was 543
$ws.Range("F43").Value = 1413   # was 1412
$ws.Range("F45").Value = 741   # was 736
$ws.Range("F46").Value = 2370   # was 2366
$ws.Range("F48").Value = 13   # was 12
$ws.Range("F49").Value = 792   # was 793
$ws.Range("F50").Value = 946   # was 945

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F22").Value = 55   # was 54

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 408   # was 407
$ws.Range("F8").Value = 5344   # was 5331
$ws.Range("F9").Value = 5344   # was 5331
$ws.Range("F10").Value = 176   # was 174
$ws.Range("F13").Value = 536   # was 535
$ws.Range("F15").Value = 1194   # was 1192
$ws.Range("F16").Value = 5291   # was 5286
$ws.Range("F19").Value = 102   # was 101
$ws.Range("F20").Value = 2950   # was 2884
$ws.Range("F24").Value = 3999   # was 3994
$ws.Range("F25").Value = 3916   # was 3912
$ws.Range("F29").Value = 254   # was 253
$ws.Range("F31").Value = 114   # was 113
$ws.Range("F32").Value = 114   # was 113
$ws.Range("F36").Value = 6963   # was 6950
$ws.Range("F37").Value = 28   # was 26
$ws.Range("F38").Value = 1141   # was 1138
$ws.Range("F39").Value = 547   # was 543
$ws.Range("F43").Value = 1413   # was 1412
$ws.Range("F45").Value = 741   # was 736
$ws.Range("F46").Value = 2370   # was 2366
$ws.Range("F48").Value = 792   # was 793
$ws.Range("F49").Value = 946   # was 945

